$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 183, 184 -- cyclic rotation of B:G by 1
$ws.Cells.Item(183,2).Value2 = 64329
$ws.Cells.Item(183,3).Value2 = "DAB-Real Activ Coconut Water Tetra 1000ml"
$ws.Cells.Item(183,4).Value2 = 120.69
$ws.Cells.Item(183,5).Value2 = 128.32
$ws.Cells.Item(183,6).Value2 = 6
$ws.Cells.Item(183,7).Value2 = 724.14
$ws.Cells.Item(184,2).Value2 = 57552
$ws.Cells.Item(184,3).Value2 = "DAB-Real Activ Coconut Water Tetra 1000ml"
$ws.Cells.Item(184,4).Value2 = 120.69
$ws.Cells.Item(184,5).Value2 = 136.86
$ws.Cells.Item(184,6).Value2 = -5
$ws.Cells.Item(184,7).Value2 = -603.45

# Rows 316, 317, 318 -- cyclic rotation of B:G by 1
$ws.Cells.Item(316,2).Value2 = 57077
$ws.Cells.Item(316,3).Value2 = "HUL-Bru Inst Poly 50g"
$ws.Cells.Item(316,4).Value2 = 93.08
$ws.Cells.Item(316,5).Value2 = 111.2
$ws.Cells.Item(316,6).Value2 = 1
$ws.Cells.Item(316,7).Value2 = 93.08
$ws.Cells.Item(317,2).Value2 = 61610
$ws.Cells.Item(317,3).Value2 = "HUL-Bru Inst Poly 50g"
$ws.Cells.Item(317,4).Value2 = 102.71
$ws.Cells.Item(317,5).Value2 = 122.71
$ws.Cells.Item(317,6).Value2 = -58
$ws.Cells.Item(317,7).Value2 = -5957.18
$ws.Cells.Item(318,2).Value2 = 63565
$ws.Cells.Item(318,3).Value2 = "HUL-Bru Inst Poly 50g"
$ws.Cells.Item(318,4).Value2 = 102.71
$ws.Cells.Item(318,5).Value2 = 109.19
$ws.Cells.Item(318,6).Value2 = 60
$ws.Cells.Item(318,7).Value2 = 6162.6

# Rows 355, 356 -- cyclic rotation of B:G by 1
$ws.Cells.Item(355,2).Value2 = 55356
$ws.Cells.Item(355,3).Value2 = "HUL-knorr schezwan 200g pch"
$ws.Cells.Item(355,4).Value2 = 47.64
$ws.Cells.Item(355,5).Value2 = 54.04
$ws.Cells.Item(355,6).Value2 = -158
$ws.Cells.Item(355,7).Value2 = -7527.12
$ws.Cells.Item(356,2).Value2 = 63510
$ws.Cells.Item(356,3).Value2 = "HUL-knorr schezwan 200g pch"
$ws.Cells.Item(356,4).Value2 = 47.64
$ws.Cells.Item(356,5).Value2 = 50.66
$ws.Cells.Item(356,6).Value2 = 167
$ws.Cells.Item(356,7).Value2 = 7955.88

# Rows 375, 376 -- cyclic rotation of B:G by 1
$ws.Cells.Item(375,2).Value2 = 63563
$ws.Cells.Item(375,3).Value2 = "HUL-lux advanced eventoned glow 4x100"
$ws.Cells.Item(375,4).Value2 = 111.96
$ws.Cells.Item(375,5).Value2 = 119.04
$ws.Cells.Item(375,6).Value2 = 15
$ws.Cells.Item(375,7).Value2 = 1679.4
$ws.Cells.Item(376,2).Value2 = 61605
$ws.Cells.Item(376,3).Value2 = "HUL-lux advanced eventoned glow 4x100"
$ws.Cells.Item(376,4).Value2 = 111.96
$ws.Cells.Item(376,5).Value2 = 133.78
$ws.Cells.Item(376,6).Value2 = -13
$ws.Cells.Item(376,7).Value2 = -1455.48

# Rows 379, 380 -- cyclic rotation of B:G by 1
$ws.Cells.Item(379,2).Value2 = 63564
$ws.Cells.Item(379,3).Value2 = "HUL-Lux Radiant Glow 3*150g"
$ws.Cells.Item(379,4).Value2 = 129.01
$ws.Cells.Item(379,5).Value2 = 137.16
$ws.Cells.Item(379,6).Value2 = 57
$ws.Cells.Item(379,7).Value2 = 7353.57
$ws.Cells.Item(380,2).Value2 = 61608
$ws.Cells.Item(380,3).Value2 = "HUL-Lux Radiant Glow 3*150g"
$ws.Cells.Item(380,4).Value2 = 129.01
$ws.Cells.Item(380,5).Value2 = 154.12
$ws.Cells.Item(380,6).Value2 = -56
$ws.Cells.Item(380,7).Value2 = -7224.56

# Rows 382, 383 -- cyclic rotation of B:G by 1
$ws.Cells.Item(382,2).Value2 = 60325
$ws.Cells.Item(382,3).Value2 = "Hul-pears pure and gentle 3x125 gm"
$ws.Cells.Item(382,4).Value2 = 126.86
$ws.Cells.Item(382,5).Value2 = 151.57
$ws.Cells.Item(382,6).Value2 = -102
$ws.Cells.Item(382,7).Value2 = -12939.72
$ws.Cells.Item(383,2).Value2 = 63560
$ws.Cells.Item(383,3).Value2 = "Hul-pears pure and gentle 3x125 gm"
$ws.Cells.Item(383,4).Value2 = 126.86
$ws.Cells.Item(383,5).Value2 = 134.87
$ws.Cells.Item(383,6).Value2 = 104
$ws.Cells.Item(383,7).Value2 = 13193.44

# Rows 431, 432 -- cyclic rotation of B:G by 1
$ws.Cells.Item(431,2).Value2 = 53082
$ws.Cells.Item(431,3).Value2 = "HUL-VIM BAR MULTIPACK FW 4X200G"
$ws.Cells.Item(431,4).Value2 = 59.47
$ws.Cells.Item(431,5).Value2 = 71.05
$ws.Cells.Item(431,6).Value2 = 1
$ws.Cells.Item(431,7).Value2 = 59.47
$ws.Cells.Item(432,2).Value2 = 63102
$ws.Cells.Item(432,3).Value2 = "HUL-Vim Bar Multipack Fw 4X200G"
$ws.Cells.Item(432,4).Value2 = 59.47
$ws.Cells.Item(432,5).Value2 = 71.05
$ws.Cells.Item(432,6).Value2 = 36
$ws.Cells.Item(432,7).Value2 = 2140.92

# Rows 457, 458 -- cyclic rotation of B:G by 1
$ws.Cells.Item(457,2).Value2 = 63681
$ws.Cells.Item(457,3).Value2 = "JLM-MBD Shiny Toothbrush Safari"
$ws.Cells.Item(457,4).Value2 = 22.42
$ws.Cells.Item(457,5).Value2 = 23.84
$ws.Cells.Item(457,6).Value2 = 65
$ws.Cells.Item(457,7).Value2 = 1457.3
$ws.Cells.Item(458,2).Value2 = 31930
$ws.Cells.Item(458,3).Value2 = "JLM-MBD Shiny Toothbrush Safari"
$ws.Cells.Item(458,4).Value2 = 22.42
$ws.Cells.Item(458,5).Value2 = 26.8
$ws.Cells.Item(458,6).Value2 = -62
$ws.Cells.Item(458,7).Value2 = -1390.04

# Rows 579, 580 -- cyclic rotation of B:G by 1
$ws.Cells.Item(579,2).Value2 = 65069
$ws.Cells.Item(579,3).Value2 = "CRE-Bourbon 100gm"
$ws.Cells.Item(579,4).Value2 = 13.45
$ws.Cells.Item(579,5).Value2 = 14.3
$ws.Cells.Item(579,6).Value2 = 172
$ws.Cells.Item(579,7).Value2 = 2313.4
$ws.Cells.Item(580,2).Value2 = 53757
$ws.Cells.Item(580,3).Value2 = "CRE-Bourbon 100gm"
$ws.Cells.Item(580,4).Value2 = 13.45
$ws.Cells.Item(580,5).Value2 = 16.08
$ws.Cells.Item(580,6).Value2 = -159
$ws.Cells.Item(580,7).Value2 = -2138.55

# Rows 583, 584 -- cyclic rotation of B:G by 1
$ws.Cells.Item(583,2).Value2 = 53263
$ws.Cells.Item(583,3).Value2 = "CRE-Butter cremfills 100gm"
$ws.Cells.Item(583,4).Value2 = 12.81
$ws.Cells.Item(583,5).Value2 = 15.29
$ws.Cells.Item(583,6).Value2 = -309
$ws.Cells.Item(583,7).Value2 = -3958.29
$ws.Cells.Item(584,2).Value2 = 65066
$ws.Cells.Item(584,3).Value2 = "CRE-Butter cremfills 100gm"
$ws.Cells.Item(584,4).Value2 = 12.81
$ws.Cells.Item(584,5).Value2 = 13.61
$ws.Cells.Item(584,6).Value2 = 313
$ws.Cells.Item(584,7).Value2 = 4009.53

# Rows 586, 587 -- cyclic rotation of B:G by 1
$ws.Cells.Item(586,2).Value2 = 64915
$ws.Cells.Item(586,3).Value2 = "CRE-Cremica Chocolate Cream 150Gm"
$ws.Cells.Item(586,4).Value2 = 19.73
$ws.Cells.Item(586,5).Value2 = 20.98
$ws.Cells.Item(586,6).Value2 = 40
$ws.Cells.Item(586,7).Value2 = 789.2
$ws.Cells.Item(587,2).Value2 = 45695
$ws.Cells.Item(587,3).Value2 = "CRE-Cremica Chocolate Cream 150Gm"
$ws.Cells.Item(587,4).Value2 = 19.73
$ws.Cells.Item(587,5).Value2 = 23.58
$ws.Cells.Item(587,6).Value2 = -36
$ws.Cells.Item(587,7).Value2 = -710.28

# Rows 590, 591 -- cyclic rotation of B:G by 1
$ws.Cells.Item(590,2).Value2 = 64922
$ws.Cells.Item(590,3).Value2 = "CRE-Cremica Golden Bytes Rich Butter 200Gm"
$ws.Cells.Item(590,4).Value2 = 19.73
$ws.Cells.Item(590,5).Value2 = 20.98
$ws.Cells.Item(590,6).Value2 = 207
$ws.Cells.Item(590,7).Value2 = 4084.11
$ws.Cells.Item(591,2).Value2 = 45706
$ws.Cells.Item(591,3).Value2 = "CRE-Cremica Golden Bytes Rich Butter 200Gm"
$ws.Cells.Item(591,4).Value2 = 19.73
$ws.Cells.Item(591,5).Value2 = 23.58
$ws.Cells.Item(591,6).Value2 = -202
$ws.Cells.Item(591,7).Value2 = -3985.46

# Rows 593, 594 -- cyclic rotation of B:G by 1
$ws.Cells.Item(593,2).Value2 = 45718
$ws.Cells.Item(593,3).Value2 = "CRE-Cremica Honey Oatmeal Cookies 50 +25 Gm"
$ws.Cells.Item(593,4).Value2 = 16.22
$ws.Cells.Item(593,5).Value2 = 19.38
$ws.Cells.Item(593,6).Value2 = -294
$ws.Cells.Item(593,7).Value2 = -4768.68
$ws.Cells.Item(594,2).Value2 = 64927
$ws.Cells.Item(594,3).Value2 = "CRE-Cremica Honey Oatmeal Cookies 50 +25 Gm"
$ws.Cells.Item(594,4).Value2 = 16.22
$ws.Cells.Item(594,5).Value2 = 17.26
$ws.Cells.Item(594,6).Value2 = 295
$ws.Cells.Item(594,7).Value2 = 4784.9

# Rows 687, 688 -- cyclic rotation of B:G by 1
$ws.Cells.Item(687,2).Value2 = 64810
$ws.Cells.Item(687,3).Value2 = "PRI-B-50 VIMAL Copper Glass 300ML (2pc Set)"
$ws.Cells.Item(687,4).Value2 = 273.92
$ws.Cells.Item(687,5).Value2 = 291.22
$ws.Cells.Item(687,6).Value2 = 7
$ws.Cells.Item(687,7).Value2 = 1917.44
$ws.Cells.Item(688,2).Value2 = 53319
$ws.Cells.Item(688,3).Value2 = "PRI-B-50 VIMAL Copper Glass 300ML (2pc Set)"
$ws.Cells.Item(688,4).Value2 = 273.92
$ws.Cells.Item(688,5).Value2 = 310.64
$ws.Cells.Item(688,6).Value2 = -6
$ws.Cells.Item(688,7).Value2 = -1643.52

# Rows 709, 710 -- cyclic rotation of B:G by 1
$ws.Cells.Item(709,2).Value2 = 64833
$ws.Cells.Item(709,3).Value2 = "Rasna 32 Glass Shikanji Nimbupani"
$ws.Cells.Item(709,4).Value2 = 32.83
$ws.Cells.Item(709,5).Value2 = 34.9
$ws.Cells.Item(709,6).Value2 = 99
$ws.Cells.Item(709,7).Value2 = 3250.17
$ws.Cells.Item(710,2).Value2 = 60025
$ws.Cells.Item(710,3).Value2 = "Rasna 32 Glass Shikanji Nimbupani"
$ws.Cells.Item(710,4).Value2 = 32.83
$ws.Cells.Item(710,5).Value2 = 37.22
$ws.Cells.Item(710,6).Value2 = -98
$ws.Cells.Item(710,7).Value2 = -3217.34

# Rows 715, 716 -- cyclic rotation of B:G by 1
$ws.Cells.Item(715,2).Value2 = 60031
$ws.Cells.Item(715,3).Value2 = "Rasna Insta Orange 500g"
$ws.Cells.Item(715,4).Value2 = 98.5
$ws.Cells.Item(715,5).Value2 = 111.69
$ws.Cells.Item(715,6).Value2 = -5
$ws.Cells.Item(715,7).Value2 = -492.5
$ws.Cells.Item(716,2).Value2 = 64836
$ws.Cells.Item(716,3).Value2 = "Rasna Insta Orange 500g"
$ws.Cells.Item(716,4).Value2 = 98.5
$ws.Cells.Item(716,5).Value2 = 104.71
$ws.Cells.Item(716,6).Value2 = 7
$ws.Cells.Item(716,7).Value2 = 689.5

# Rows 720, 721 -- cyclic rotation of B:G by 1
$ws.Cells.Item(720,2).Value2 = 64830
$ws.Cells.Item(720,3).Value2 = "Rasna Nagpur Orange (32 Glass)"
$ws.Cells.Item(720,4).Value2 = 32.83
$ws.Cells.Item(720,5).Value2 = 34.9
$ws.Cells.Item(720,6).Value2 = 117
$ws.Cells.Item(720,7).Value2 = 3841.11
$ws.Cells.Item(721,2).Value2 = 60022
$ws.Cells.Item(721,3).Value2 = "Rasna Nagpur Orange (32 Glass)"
$ws.Cells.Item(721,4).Value2 = 32.83
$ws.Cells.Item(721,5).Value2 = 37.22
$ws.Cells.Item(721,6).Value2 = -113
$ws.Cells.Item(721,7).Value2 = -3709.79

# Rows 872, 873 -- cyclic rotation of B:G by 1
$ws.Cells.Item(872,2).Value2 = 65079
$ws.Cells.Item(872,3).Value2 = "Shankys Tip Top Hing Jeera Peanut/ Salted Peanut 200 Gm"
$ws.Cells.Item(872,4).Value2 = 40.87
$ws.Cells.Item(872,5).Value2 = 43.44
$ws.Cells.Item(872,6).Value2 = 21
$ws.Cells.Item(872,7).Value2 = 858.27
$ws.Cells.Item(873,2).Value2 = 54751
$ws.Cells.Item(873,3).Value2 = "Shankys Tip Top Hing Jeera Peanut/ Salted Peanut 200 Gm"
$ws.Cells.Item(873,4).Value2 = 40.87
$ws.Cells.Item(873,5).Value2 = 46.34
$ws.Cells.Item(873,6).Value2 = -19
$ws.Cells.Item(873,7).Value2 = -776.53
